# Automatische test-sync: 2025-06-19 17:47:21
# Append a new "Afmelding nieuwsbrief" row to the Logs sheet and
# refresh the dependent conditional formatting ranges and the
# Dashboard summary count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Add the new log entry in row 19 ---------------------------------
$ws.Range("A19").Value = "Afmelding nieuwsbrief"
$ws.Range("B19").Value = "mailmind.test@zohomail.eu"
$ws.Range("C19").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Range("D19").Value = "Afmelding"
$ws.Range("F19").Value = "2025-06-19 16:58:11"
$ws.Range("G19").Value = "Nee"

# --- Extend the conditional formatting ranges to include row 19 ------
$catRules = $ws.Range("D2:D18").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($ws.Range("D2:D19"))
}

$answeredRules = $ws.Range("G2:G18").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($ws.Range("G2:G19"))
}

# --- Update the Dashboard summary count for "Afmelding" --------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 4
